$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order / labels (row 1)
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"

# New data values (rows 2-7)
$data = @(
    @(0,1,0,0,0,0),
    @(0,0,0,0,1,0),
    @(0,0,0,1,0,0),
    @(0,0,0,0,0,1),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}
